# Daily attendance processing - 2025-10-10 19:40:45
# Swap the order of the two comma-separated "Recorded By" names in column G
# for the specified rows (the underlying data did not change, only the
# display order of the two collaborators listed in each cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(3,6,7,10,11,12,13,14,15,30,33,34,37,38,39,40,41,42,57,60,61,64,65,66,67,68,69,86,87,88,89,90,93,95,112,113,114,115,116,119,121,138,139,140,141,142,145,147)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value()
    if ($val -ne $null) {
        $parts = $val -split ",\s*"
        if ($parts.Count -eq 2) {
            $cell.Value = ($parts[1].Trim() + ", " + $parts[0].Trim())
        }
    }
}
